$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '99.002.38'
$ws.Range("E2").Value = '  +0.86%  '
# Row 3
$ws.Range("D3").Value = '3.313.12'
$ws.Range("E3").Value = '  -1.64%  '
# Row 4
$ws.Range("E4").Value = '  -0.05%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '255.88'
$ws.Range("E5").Value = '  +0.56%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '626.20'
$ws.Range("E6").Value = '  +0.54%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.46'
$ws.Range("E7").Value = '  +21.38%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.417'
$ws.Range("E8").Value = '  +8.15%  '
# Row 9
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.04'
$ws.Range("E9").Value = '  +28.37%  '
# Row 10
$ws.Range("B10").Value = 'USDC'
$ws.Range("C10").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.999'
$ws.Range("E10").Value = '  -0.06%  '
# Row 11
$ws.Range("D11").Value = '3.311.41'
$ws.Range("E11").Value = '  -1.58%  '
# Row 12
$ws.Range("E12").Value = '  +2.84%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.45'
$ws.Range("E13").Value = '  +15.65%  '
# Row 14
$ws.Range("D14").Value = '98.632.94'
$ws.Range("E14").Value = '  +0.78%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000254'
$ws.Range("E15").Value = '  +2.85%  '
# Row 16
$ws.Range("D16").Value = '3.936.93'
$ws.Range("E16").Value = '  -1.37%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.40'
$ws.Range("E17").Value = '  -1.60%  '
# Row 18
$ws.Range("D18").Value = '3.313.59'
$ws.Range("E18").Value = '  -1.57%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.47'
$ws.Range("E19").Value = '  -4.85%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.77'
$ws.Range("E20").Value = '  +6.62%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.42'
$ws.Range("E21").Value = '  +9.03%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '488.42'
$ws.Range("E22").Value = '  +1.45%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.48'
$ws.Range("E23").Value = '  +3.35%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000202'
$ws.Range("E24").Value = '  -2.64%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.88'
$ws.Range("E25").Value = '  +3.08%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.342'
$ws.Range("E26").Value = '  +35.76%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '89.34'
$ws.Range("E27").Value = '  +1.43%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.19'
$ws.Range("E28").Value = '  +1.39%  '
# Row 29
$ws.Range("D29").Value = '3.486.29'
$ws.Range("E29").Value = '  -1.74%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.149'
$ws.Range("E30").Value = '  +17.99%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.12%  '
# Row 32
$ws.Range("E32").Value = '  +1.96%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.64'
$ws.Range("E33").Value = '  +14.82%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.20%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '28.11'
$ws.Range("E35").Value = '  +2.70%  '
# Row 36
$ws.Range("B36").Value = 'PolygonEcosystemToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.486'
$ws.Range("E36").Value = '  +8.30%  '
# Row 37
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.152'
$ws.Range("E37").Value = '  -0.41%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.36'
$ws.Range("E38").Value = '  +0.60%  '
# Row 39
$ws.Range("E39").Value = '  +1.25%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '497.46'
$ws.Range("E40").Value = '  -5.17%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '24.74'
$ws.Range("E41").Value = '  -0.32%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.88'
$ws.Range("E42").Value = '  +1.99%  '
# Row 43
$ws.Range("E43").Value = '  -0.80%  '
# Row 44
$ws.Range("E44").Value = '  +0.13%  '
# Row 45
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  -0.01%  '
# Row 46
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.17'
$ws.Range("E46").Value = '  -2.48%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '160.70'
$ws.Range("E47").Value = '  -0.30%  '
# Row 48
$ws.Range("E48").Value = '  +1.83%  '
# Row 49
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.861'
$ws.Range("E49").Value = '  +8.02%  '
# Row 50
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.78'
$ws.Range("E50").Value = '  +5.19%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.34'
$ws.Range("E51").Value = '  +14.41%  '
